# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Sat Apr  6 22:23:46 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are stored as text in the workbook (e.g. "68.366.82",
# "1.00", "0.120"); Excel would otherwise silently reinterpret some of them as
# numbers (dropping trailing zeros / turning "1.00" into 1). Force text storage,
# then restore the default "Normal" style so no visible formatting changes.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "68.366.82"
Set-TextValue $ws.Range("D3") "3.352.26"
$ws.Range("E3").Value = "  +0.87%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "584.89"
$ws.Range("E5").Value = "  +0.89%  "
Set-TextValue $ws.Range("D6") "177.34"
$ws.Range("E6").Value = "  +1.59%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  +4.01%  "
Set-TextValue $ws.Range("D10") "0.582"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("E11").Value = "  +6.04%  "
Set-TextValue $ws.Range("D12") "0.0000274"
$ws.Range("E12").Value = "  +2.04%  "
Set-TextValue $ws.Range("D13") "691.84"
$ws.Range("E13").Value = "  +4.51%  "
Set-TextValue $ws.Range("D14") "3.904.48"
$ws.Range("E14").Value = "  +0.97%  "
Set-TextValue $ws.Range("D15") "8.47"
$ws.Range("E15").Value = "  +0.93%  "
Set-TextValue $ws.Range("D16") "68.419.48"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D17") "0.120"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D18") "3.340.16"
$ws.Range("E18").Value = "  +0.39%  "
Set-TextValue $ws.Range("D19") "17.50"
$ws.Range("E19").Value = "  +1.03%  "
Set-TextValue $ws.Range("D20") "11.23"
$ws.Range("E20").Value = "  +2.79%  "
Set-TextValue $ws.Range("D21") "0.895"
$ws.Range("E21").Value = "  +1.10%  "
Set-TextValue $ws.Range("D22") "5.49"
$ws.Range("E22").Value = "  +2.50%  "
Set-TextValue $ws.Range("D23") "16.95"
$ws.Range("E23").Value = "  -0.44%  "
Set-TextValue $ws.Range("D24") "100.04"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("E25").Value = "  +1.89%  "
Set-TextValue $ws.Range("D26") "2.71"
$ws.Range("E26").Value = "  +1.88%  "
Set-TextValue $ws.Range("D27") "9.52"
$ws.Range("E27").Value = "  +3.23%  "
Set-TextValue $ws.Range("D28") "33.05"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("E29").Value = "  +1.64%  "
Set-TextValue $ws.Range("D30") "6.98"
$ws.Range("E30").Value = "  -4.08%  "
Set-TextValue $ws.Range("D32") "552.19"
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("E33").Value = "  +0.70%  "
Set-TextValue $ws.Range("D34") "58.14"
$ws.Range("E34").Value = "  +2.82%  "
Set-TextValue $ws.Range("D35") "3.717.02"
$ws.Range("E35").Value = "  +1.18%  "
$ws.Range("E36").Value = "  -0.07%  "
Set-TextValue $ws.Range("D37") "3.41"
$ws.Range("E37").Value = "  +4.47%  "
Set-TextValue $ws.Range("D38") "0.141"
$ws.Range("E38").Value = "  +9.23%  "
Set-TextValue $ws.Range("D39") "34.73"
$ws.Range("E39").Value = "  +1.75%  "
Set-TextValue $ws.Range("D40") "3.19"
$ws.Range("E40").Value = "  +2.88%  "
Set-TextValue $ws.Range("D41") "2.62"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("E43").Value = "  +0.98%  "
Set-TextValue $ws.Range("D44") "3.25"
$ws.Range("E44").Value = "  -2.90%  "
Set-TextValue $ws.Range("D45") "0.0413"
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("E48").Value = "  +0.02%  "
Set-TextValue $ws.Range("D49") "1.34"
$ws.Range("E49").Value = "  -1.44%  "
Set-TextValue $ws.Range("D50") "131.84"
$ws.Range("E50").Value = "  +1.74%  "
Set-TextValue $ws.Range("D51") "2.62"
$ws.Range("E51").Value = "  -0.62%  "
